$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91..124 down to 92..125
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record
$ws.Range("A91").Value = 5
$ws.Range("B91").Value = "Macroferia Regional de Talca"
$ws.Range("C91").Value = "Maule"
$ws.Range("D91").Value = 45229
$ws.Range("E91").Value = 7
$ws.Range("F91").Value = 300000000
$ws.Range("G91").Value = "Espárragos"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 1100
$ws.Range("L91").Value = 1100
$ws.Range("M91").Value = 1100
$ws.Range("N91").Value = "`$/kilo"
$ws.Range("O91").Value = "Provincia de Linares"
$ws.Range("P91").Value = 1100
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
